$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E columns to remain text (values are numeric-looking strings, e.g. "1.00", "0.0241")
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "62.815.79"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "2.681.92"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "556.15"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "157.67"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("D12").Value = "5.39"
$ws.Range("E12").Value = "  -3.72%  "
$ws.Range("D13").Value = "3.156.60"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "26.48"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "62.768.14"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "0.0000146"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("D17").Value = "2.685.97"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "11.85"
$ws.Range("E18").Value = "  -4.20%  "
$ws.Range("D19").Value = "4.61"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("D20").Value = "345.01"
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("D21").Value = "6.21"
$ws.Range("E21").Value = "  -5.35%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "0.505"
$ws.Range("E23").Value = "  -3.03%  "
$ws.Range("D24").Value = "63.15"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "8.15"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0854"
$ws.Range("E28").Value = "  -5.31%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "1.39"
$ws.Range("E29").Value = "  +5.51%  "
$ws.Range("D30").Value = "7.23"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "164.03"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.47"
$ws.Range("D36").Value = "19.46"
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").Value = "346.31"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "6.19"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "0.937"
$ws.Range("E40").Value = "  -3.79%  "
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("D42").Value = "38.35"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "20.77"
$ws.Range("E43").Value = "  -4.73%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "20.15"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0556"
$ws.Range("E47").Value = "  -4.28%  "
$ws.Range("D48").Value = "10.99"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "0.0969"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "128.70"
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0241"
$ws.Range("E51").Value = "  -3.37%  "
